# إضافة حدث جديد في Card23 by admin at 2025-12-08 07:42:00
#
# Fills in the previously-blank B18:K18 cells with the literal text "nan"
# (matching the rest of the sheet's "missing value" convention) and appends
# a brand-new service-log row 19 for card 23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Card23")

# --- Row 18: the blank tracking cells (B..K) were left truly empty; backfill
#     them with the sheet's usual "nan" placeholder text. ---
$row18Cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
foreach ($col in $row18Cols) {
    $ws.Range($col + "18").Value = "nan"
}

# --- Row 19: new service event for card 23. ---
# A19 holds the numeric-looking card id "23" as *text* (matching the rest of
# column A), so force a text format before writing it, then drop the format
# again so no stray style sticks to the cell.
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "23"
$ws.Range("A19").ClearFormats()

# B19..K19 stay blank, same as the rest of the new row.

$ws.Range("L19").Value = "14\8\2025"
$ws.Range("M19").Value = "804 t"
$ws.Range("N19").Value = "تم تغير زيت الجيربوكس"
$ws.Range("O19").Value = "تم العمل"
